$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the three new vocabulary rows (English / Chinese / Vietnamese) ---
$ws.Range("A56").Value = "cheeze"
$ws.Range("B56").Value = "起司"
$ws.Range("C56").Value = "phô mai"

$ws.Range("A57").Value = "salt"
$ws.Range("B57").Value = "鹽"
$ws.Range("C57").Value = "muối"

# --- Refresh the AutoFilter range (header + data through row 57) before
#     row 58 lands, so the filter ref stays A1:C57 like the saved file. ---
$ws.Range("A1:C57").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=工作表1!`$A`$1:`$C`$57"
    }
}

$ws.Range("A58").Value = "lemon juice"
$ws.Range("B58").Value = "檸檬汁"
$ws.Range("C58").Value = "nước chanh"

# Row 58 mirrors the wrapped/top-aligned header style used by column A on
# other "taller" rows (e.g. A2, A14, A20) and gets the same custom height.
$ws.Range("A58").WrapText = $true
$ws.Range("A58").HorizontalAlignment = -4131
$ws.Range("A58").VerticalAlignment = -4160
$ws.Range("A58").Font.Name = "Microsoft JhengHei Light"
$ws.Range("A58").Font.Size = 12
$ws.Range("A58:C58").RowHeight = 19.2

# --- View: zoom + selection matching the saved workbook state ---
$ws.Application.ActiveWindow.Zoom = 79
$ws.Range("B59").Select()
